$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2081911262798635
$ws.Range("C2").Value = 0.552901023890785
$ws.Range("J2").Value = 0.006825938566552901
$ws.Range("P2").Value = 0.1331058020477816
$ws.Range("S2").Value = 0.09897610921501707
$ws.Range("J3").Value = 0.03105590062111801
$ws.Range("P3").Value = 0.7391304347826086
$ws.Range("S3").Value = 0.2298136645962733
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.7105263157894737
$ws.Range("S4").Value = 0.2368421052631579
$ws.Range("B6").Value = 0.08239700374531835
$ws.Range("D6").Value = 0.01123595505617977
$ws.Range("F6").Value = 0.0449438202247191
$ws.Range("J6").Value = 0.2059925093632959
$ws.Range("O6").Value = 0.01872659176029963
$ws.Range("Q6").Value = 0.1797752808988764
$ws.Range("R6").Value = 0.0749063670411985
$ws.Range("S6").Value = 0.3820224719101123
$ws.Range("B7").Value = 0.1065573770491803
$ws.Range("D7").Value = 0.01639344262295082
$ws.Range("F7").Value = 0.04098360655737705
$ws.Range("J7").Value = 0.1229508196721311
$ws.Range("O7").Value = 0.02459016393442623
$ws.Range("Q7").Value = 0.180327868852459
$ws.Range("R7").Value = 0.08196721311475409
$ws.Range("S7").Value = 0.4262295081967213
$ws.Range("B8").Value = 0.08607594936708861
$ws.Range("D8").Value = 0.02278481012658228
$ws.Range("F8").Value = 0.07088607594936709
$ws.Range("J8").Value = 0.1139240506329114
$ws.Range("O8").Value = 0.01518987341772152
$ws.Range("Q8").Value = 0.1721518987341772
$ws.Range("R8").Value = 0.0759493670886076
$ws.Range("S8").Value = 0.4430379746835443
$ws.Range("B9").Value = 0.09322033898305085
$ws.Range("D9").Value = 0.0211864406779661
$ws.Range("F9").Value = 0.1059322033898305
$ws.Range("J9").Value = 0.1016949152542373
$ws.Range("O9").Value = 0.0211864406779661
$ws.Range("Q9").Value = 0.1779661016949153
$ws.Range("R9").Value = 0.07203389830508475
$ws.Range("S9").Value = 0.4067796610169492
$ws.Range("B10").Value = 0.1
$ws.Range("D10").Value = 0.01438848920863309
$ws.Range("E10").Value = 0.001438848920863309
$ws.Range("F10").Value = 0.08345323741007195
$ws.Range("J10").Value = 0.102158273381295
$ws.Range("O10").Value = 0.01942446043165467
$ws.Range("Q10").Value = 0.1971223021582734
$ws.Range("R10").Value = 0.08633093525179857
$ws.Range("S10").Value = 0.3956834532374101
$ws.Range("G11").Value = 0.1368821292775665
$ws.Range("J11").Value = 0.1330798479087452
$ws.Range("K11").Value = 0.2395437262357414
$ws.Range("L11").Value = 0.4638783269961977
$ws.Range("S11").Value = 0.02661596958174905
$ws.Range("G12").Value = 0.5581395348837209
$ws.Range("J12").Value = 0.3023255813953488
$ws.Range("L12").Value = 0.04651162790697674
$ws.Range("S12").Value = 0.09302325581395349
$ws.Range("G13").Value = 0.4705882352941176
$ws.Range("J13").Value = 0.4117647058823529
$ws.Range("S13").Value = 0.1176470588235294
$ws.Range("F15").Value = 0.02222222222222222
$ws.Range("H15").Value = 0.1185185185185185
$ws.Range("J15").Value = 0.4370370370370371
$ws.Range("K15").Value = 0.05185185185185185
$ws.Range("M15").Value = 0.01481481481481482
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.2148148148148148
$ws.Range("F16").Value = 0.01694915254237288
$ws.Range("H16").Value = 0.096045197740113
$ws.Range("I16").Value = 0.1016949152542373
$ws.Range("J16").Value = 0.4519774011299435
$ws.Range("K16").Value = 0.1186440677966102
$ws.Range("M16").Value = 0.005649717514124294
$ws.Range("N16").Value = 0.005649717514124294
$ws.Range("O16").Value = 0.07909604519774012
$ws.Range("S16").Value = 0.1242937853107345
$ws.Range("F17").Value = 0.03083700440528634
$ws.Range("H17").Value = 0.1475770925110132
$ws.Range("I17").Value = 0.09251101321585903
$ws.Range("J17").Value = 0.4449339207048458
$ws.Range("K17").Value = 0.05726872246696035
$ws.Range("M17").Value = 0.01541850220264317
$ws.Range("O17").Value = 0.09911894273127753
$ws.Range("S17").Value = 0.1123348017621145
$ws.Range("H18").Value = 0.1658291457286432
$ws.Range("I18").Value = 0.09547738693467336
$ws.Range("J18").Value = 0.3919597989949749
$ws.Range("K18").Value = 0.1055276381909548
$ws.Range("M18").Value = 0.01507537688442211
$ws.Range("O18").Value = 0.08040201005025126
$ws.Range("S18").Value = 0.1457286432160804
$ws.Range("F19").Value = 0.01621223286661754
$ws.Range("H19").Value = 0.180545320560059
$ws.Range("I19").Value = 0.1016949152542373
$ws.Range("J19").Value = 0.4089904200442152
$ws.Range("K19").Value = 0.08769344141488578
$ws.Range("M19").Value = 0.0154753131908622
$ws.Range("N19").Value = 0.0007369196757553427
$ws.Range("O19").Value = 0.08327192336035372
$ws.Range("S19").Value = 0.105379513633014
